# Update products to account for dairy:
#  - Several rows get track_inventory (Q) / visible (S) flipped from
#    "True" to "False", and their stock_inventory (R) zeroed out.
#  - Several other rows just get their stock_inventory (R) count reduced.
#  - Two placeholder dairy product IDs (B) are replaced with real IDs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose track_inventory (Q) and visible (S) flip from True -> False,
# with stock_inventory (R) reset to 0.
$falseRows = @(12, 22, 27, 35, 37, 39, 41, 44, 47, 48, 49, 69, 115)

# Use a cell that already holds the text "False" (shared string) as the
# copy source so the pasted cells stay text cells instead of being
# re-parsed into native booleans.
$ws.Range("O142").Copy() | Out-Null
foreach ($r in $falseRows) {
    $ws.Range("Q$r").PasteSpecial(-4163) | Out-Null
    $ws.Range("S$r").PasteSpecial(-4163) | Out-Null
}
$ws.Application.CutCopyMode = 0

# stock_inventory (R) updates for the rows above (now 0)
$ws.Range("R12").Value = 0
$ws.Range("R22").Value = 0
$ws.Range("R27").Value = 0
$ws.Range("R35").Value = 0
$ws.Range("R37").Value = 0
$ws.Range("R39").Value = 0
$ws.Range("R41").Value = 0
$ws.Range("R44").Value = 0
$ws.Range("R47").Value = 0
$ws.Range("R48").Value = 0
$ws.Range("R49").Value = 0
$ws.Range("R69").Value = 0
$ws.Range("R115").Value = 0

# stock_inventory (R) adjustments for rows whose track_inventory/visible
# flags are unchanged.
$ws.Range("R13").Value = 13
$ws.Range("R19").Value = 4
$ws.Range("R20").Value = 4
$ws.Range("R26").Value = 8
$ws.Range("R32").Value = 3
$ws.Range("R33").Value = 2
$ws.Range("R36").Value = 4
$ws.Range("R38").Value = 8
$ws.Range("R46").Value = 20
$ws.Range("R50").Value = 4
$ws.Range("R72").Value = 5

# localLineProductID (B) updates - placeholder dairy SKUs replaced with
# real ones.
$ws.Range("B142").Value = 990715
$ws.Range("B155").Value = 990712
